$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param(
        [string]$Addr,
        [string]$Text
    )
    $c = $ws.Range($Addr)
    $c.NumberFormat = "@"
    $c.Value = $Text
    $c.Style = "Normal"
}

Set-CellText "D2" "33.733.43"
Set-CellText "E2" "  -0.25%  "
Set-CellText "D3" "1.764.75"
Set-CellText "E3" "  -0.67%  "
Set-CellText "E4" "  +0.28%  "
Set-CellText "D5" "224.49"
Set-CellText "E5" "  +1.70%  "
Set-CellText "D6" "0.544"
Set-CellText "E6" "  -1.21%  "
Set-CellText "E7" "  +0.22%  "
Set-CellText "D8" "31.99"
Set-CellText "E8" "  +3.06%  "
Set-CellText "E9" "  +0.75%  "
Set-CellText "E10" "  -3.06%  "
Set-CellText "E11" "  +1.68%  "
Set-CellText "D12" "2.019.99"
Set-CellText "E12" "  -0.51%  "
Set-CellText "D13" "11.21"
Set-CellText "E13" "  +6.72%  "
Set-CellText "D14" "1.765.49"
Set-CellText "E14" "  -0.53%  "
Set-CellText "D15" "33.728.55"
Set-CellText "E15" "  -0.27%  "
Set-CellText "D16" "0.610"
Set-CellText "E16" "  -2.39%  "
Set-CellText "E17" "  -1.79%  "
Set-CellText "D18" "66.54"
Set-CellText "E18" "  -1.82%  "
Set-CellText "B19" "ShibaInu"
Set-CellText "C19" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-CellText "D19" "0.0₃0771"
Set-CellText "E19" "  -0.22%  "
Set-CellText "B20" "BitcoinCash"
Set-CellText "C20" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-CellText "D20" "237.56"
Set-CellText "E20" "  -2.73%  "
Set-CellText "E21" "  +0.23%  "
Set-CellText "D22" "10.57"
Set-CellText "D23" "4.04"
Set-CellText "E23" "  -0.75%  "
Set-CellText "E24" "  -1.74%  "
Set-CellText "D25" "159.08"
Set-CellText "E25" "  +1.22%  "
Set-CellText "D26" "16.11"
Set-CellText "E26" "  -1.49%  "
Set-CellText "D27" "7.02"
Set-CellText "E27" "  +0.74%  "
Set-CellText "E28" "  -0.18%  "
Set-CellText "E29" "  +0.37%  "
Set-CellText "D30" "1.23"
Set-CellText "D31" "0.0509"
Set-CellText "E31" "  -2.15%  "
Set-CellText "D32" "3.58"
Set-CellText "E32" "  -3.05%  "
Set-CellText "D33" "3.50"
Set-CellText "E33" "  +0.36%  "
Set-CellText "D34" "1.78"
Set-CellText "E34" "  -1.05%  "
Set-CellText "D35" "1.379.98"
Set-CellText "E35" "  -0.92%  "
Set-CellText "D36" "0.653"
Set-CellText "E36" "  +2.55%  "
Set-CellText "E37" "  -1.36%  "
Set-CellText "E38" "  -0.25%  "
Set-CellText "E39" "  +6.03%  "
Set-CellText "E40" "  +0.66%  "
Set-CellText "B41" "ARBITRUM"
Set-CellText "C41" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-CellText "D41" "0.908"
Set-CellText "E41" "  -2.30%  "
Set-CellText "B42" "InjectiveProtocol"
Set-CellText "C42" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-CellText "D42" "13.61"
Set-CellText "E42" "  +16.44%  "
Set-CellText "D43" "77.58"
Set-CellText "E43" "  -1.54%  "
Set-CellText "E44" "  -1.66%  "
Set-CellText "E45" "  +4.30%  "
Set-CellText "D46" "0.0₆0138"
Set-CellText "E46" "  +15.41%  "
Set-CellText "E47" "  +1.90%  "
Set-CellText "D48" "107.44"
Set-CellText "E48" "  +2.98%  "
Set-CellText "D49" "5.81"
Set-CellText "E49" "  -0.95%  "
Set-CellText "D50" "1.920.50"
Set-CellText "E50" "  +0.49%  "
Set-CellText "E51" "  +0.51%  "

Write-Host "Applied 88 changes"
